# Fruta / hortaliza, semanal
# Insert a new weekly record at row 461 for "Feria Lagunitas de Puerto Montt"
# (Naranja / Fukumoto / Segunda), pushing the existing rows 461:498 down to
# 462:499.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 461, shifting rows 461-498
# down to 462-499 (Excel's normal EntireRow.Insert semantics).
$ws.Rows.Item(461).Insert()

# Populate the newly inserted row 461 with the new observation.
$ws.Range("A461").Value2 = 4
$ws.Range("B461").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C461").Value2 = "Los Lagos"
$ws.Range("D461").Value2 = 44769
$ws.Range("E461").Value2 = 10
$ws.Range("F461").Value2 = "Fruta"
$ws.Range("G461").Value2 = 100102
$ws.Range("H461").Value2 = "Cítricos"
$ws.Range("I461").Value2 = 100102005
$ws.Range("J461").Value2 = "Naranja"
$ws.Range("K461").Value2 = "Fukumoto"
$ws.Range("L461").Value2 = "Segunda"
$ws.Range("M461").Value2 = 200
$ws.Range("N461").Value2 = 8000
$ws.Range("O461").Value2 = 8500
$ws.Range("P461").Value2 = 8250
$ws.Range("Q461").Value2 = "$/malla 16 kilos"
$ws.Range("R461").Value2 = "Región de O'Higgins"
$ws.Range("S461").Value2 = 516
$ws.Range("T461").Value2 = 16
